$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '''308.87'
$ws.Range("E2").Value = '''-3.82%'
$ws.Range("G2").Value = '''17'

$ws.Range("D3").Value = '''54.51'
$ws.Range("E3").Value = '''11.00%'
$ws.Range("G3").Value = '''17'

$ws.Range("D4").Value = '''5.104'
$ws.Range("E4").Value = '''-4.35%'
$ws.Range("G4").Value = '''17'

$ws.Range("D5").Value = '''0.07861'
$ws.Range("E5").Value = '''-2.40%'
$ws.Range("G5").Value = '''17'

$ws.Range("D6").Value = '''4.554'
$ws.Range("E6").Value = '''-1.26%'
$ws.Range("G6").Value = '''17'

$ws.Range("D7").Value = '''1.383'
$ws.Range("E7").Value = '''-1.13%'
$ws.Range("G7").Value = '''17'

$ws.Range("D8").Value = '''1.731'
$ws.Range("E8").Value = '''5.50%'
$ws.Range("G8").Value = '''17'

$ws.Range("D9").Value = '''0.1242'
$ws.Range("G9").Value = '''17'

$ws.Range("D10").Value = '''0.2007'
$ws.Range("E10").Value = '''1.59%'
$ws.Range("G10").Value = '''17'

$ws.Range("D11").Value = '''0.04722'
$ws.Range("E11").Value = '''0.72%'
$ws.Range("G11").Value = '''17'

$ws.Range("D12").Value = '''0.09422'
$ws.Range("E12").Value = '''-1.66%'
$ws.Range("G12").Value = '''17'

$ws.Range("D13").Value = '''0.1043'
$ws.Range("E13").Value = '''-0.25%'
$ws.Range("G13").Value = '''17'

$ws.Range("D14").Value = '''0.001264'
$ws.Range("E14").Value = '''-4.64%'
$ws.Range("G14").Value = '''17'

$ws.Range("D15").Value = '''0.005655'
$ws.Range("E15").Value = '''-3.23%'
$ws.Range("G15").Value = '''17'

$ws.Range("E16").Value = '''2,015.92%'
$ws.Range("G16").Value = '''17'

$ws.Range("B17").Value = 'HotbitToken'
$ws.Range("C17").Value = 'https://coinranking.com/coin/uQJB8Ocu8lTb+hotbittoken-htb'
$ws.Range("D17").Value = '''0.003943'
$ws.Range("E17").Value = '''-8.46%'
$ws.Range("G17").Value = '''17'

$ws.Range("B18").Value = 'LEO'
$ws.Range("C18").Value = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
$ws.Range("D18").Value = '''3.328'
$ws.Range("E18").Value = '''-0.44%'
$ws.Range("G18").Value = '''17'

$ws.Range("B19").Value = 'BTSEToken'
$ws.Range("C19").Value = 'https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse'
$ws.Range("D19").Value = '''2.413'
$ws.Range("E19").Value = '''-1.42%'
$ws.Range("G19").Value = '''17'

$ws.Range("B20").Value = 'BitpandaEcosystemToken'
$ws.Range("C20").Value = 'https://coinranking.com/coin/Uzf_Wjqc+bitpandaecosystemtoken-best'
$ws.Range("D20").Value = '''0.3419'
$ws.Range("E20").Value = '''-2.56%'
$ws.Range("G20").Value = '''17'

$ws.Range("B21").Value = 'MCDex'
$ws.Range("C21").Value = 'https://coinranking.com/coin/3nMM61qeg+mcdex-mcb'
$ws.Range("D21").Value = '''8.357'
$ws.Range("E21").Value = '''4.27%'
$ws.Range("G21").Value = '''17'

$ws.Range("B22").Value = 'ProBitToken'
$ws.Range("C22").Value = 'https://coinranking.com/coin/lQP4d6T2+probittoken-prob'
$ws.Range("D22").Value = '''0.1362'
$ws.Range("E22").Value = '''-0.85%'
$ws.Range("G22").Value = '''17'

$ws.Range("B23").Value = 'ZBToken'
$ws.Range("C23").Value = 'https://coinranking.com/coin/CxmvOsCyENPso+zbtoken-zb'
$ws.Range("D23").Value = '''0.2915'
$ws.Range("E23").Value = '''-5.76%'
$ws.Range("G23").Value = '''17'

$ws.Range("B24").Value = 'CoinExToken'
$ws.Range("C24").Value = 'https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet'
$ws.Range("D24").Value = '''0.04158'
$ws.Range("E24").Value = '''-0.82%'
$ws.Range("G24").Value = '''17'

$ws.Range("B25").Value = 'BitKan'
$ws.Range("C25").Value = 'https://coinranking.com/coin/RDOsLDgvY-AXe+bitkan-kan'
$ws.Range("D25").Value = '''0.001261'
$ws.Range("E25").Value = '''-3.99%'
$ws.Range("G25").Value = '''17'

$ws.Range("E26").Value = '''-0.01%'
$ws.Range("G26").Value = '''17'

$ws.Range("G27").Value = '''17'

$ws.Range("G28").Value = '''17'

$ws.Range("G29").Value = '''17'

$ws.Range("G30").Value = '''17'

$ws.Range("G31").Value = '''17'

$ws.Range("G32").Value = '''17'

$ws.Range("G33").Value = '''17'

$ws.Range("G34").Value = '''17'

$ws.Range("G35").Value = '''17'

$ws.Range("G36").Value = '''17'

$ws.Range("G37").Value = '''17'

$ws.Range("D38").Value = '''0.02599'
$ws.Range("E38").Value = '''-5.04%'
$ws.Range("G38").Value = '''17'

$ws.Range("D39").Value = '''0.05858'
$ws.Range("E39").Value = '''-6.76%'
$ws.Range("G39").Value = '''17'

$ws.Range("D40").Value = '''0.01074'
$ws.Range("E40").Value = '''-1.11%'
$ws.Range("G40").Value = '''17'

$ws.Range("D41").Value = '''0.007953'
$ws.Range("E41").Value = '''-0.95%'
$ws.Range("G41").Value = '''17'

$ws.Range("D42").Value = '''0.1372'
$ws.Range("E42").Value = '''-6.31%'
$ws.Range("G42").Value = '''17'

$ws.Range("D43").Value = '''0.008217'
$ws.Range("E43").Value = '''4.08%'
$ws.Range("G43").Value = '''17'

$ws.Range("D44").Value = '''0.008366'
$ws.Range("E44").Value = '''-3.48%'
$ws.Range("G44").Value = '''17'

$ws.Range("D45").Value = '''0.3371'
$ws.Range("E45").Value = '''-3.76%'
$ws.Range("G45").Value = '''17'

$ws.Range("D46").Value = '''0.00007313'
$ws.Range("E46").Value = '''6.75%'
$ws.Range("G46").Value = '''17'

$ws.Range("D47").Value = '''0.00000000750'
$ws.Range("E47").Value = '''-0.04%'
$ws.Range("G47").Value = '''17'

$ws.Range("D48").Value = '''0.05687'
$ws.Range("E48").Value = '''-4.46%'
$ws.Range("G48").Value = '''17'

$ws.Range("E49").Value = '''-34.56%'
$ws.Range("G49").Value = '''17'

$ws.Range("D50").Value = '''0.00002099'
$ws.Range("E50").Value = '''-0.04%'
$ws.Range("G50").Value = '''17'

$ws.Range("D51").Value = '''0.0001999'
$ws.Range("E51").Value = '''-0.04%'
$ws.Range("G51").Value = '''17'
